$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update formula in C2: drop the RIGHT(...) trimming, just concatenate first initial with full last name
$ws.Range("C2").Formula = "=LEFT(A2,1)&B2"

# Update the active selection on the sheet to C3
$ws.Range("C3").Select()
